$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column J holds the 2022 figures, mirroring the formatting of column I.
$ws.Range("J4").Value = 2022
$ws.Range("J5").Value = 96.4
$ws.Range("J6").Value = 96.4
$ws.Range("J7").Value = 97.9
$ws.Range("J8").Value = 95.3
$ws.Range("J9").Value = 93.8
$ws.Range("J10").Value = 95.5
$ws.Range("J11").Value = 94.4
$ws.Range("J12").Value = 95
$ws.Range("J13").Value = 98.7
$ws.Range("J14").Value = 97.3

# Copy column I's formatting (including the thin border line in row 3/4/14,
# number formats, fonts, etc.) into the new column J.
$ws.Range("I3:I14").Copy()
$ws.Range("J3:J14").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Restore the selection state recorded in the saved workbook.
$ws.Range("L10").Select()
